$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(76, 8).Value = 3101.4546
$ws.Cells.Item(76, 9).Value = 3030.6924
$ws.Cells.Item(76, 10).Value = 3364.2856
$ws.Cells.Item(76, 11).Value = 3030.6924
$ws.Cells.Item(76, 12).Value = 3364.2856
$ws.Cells.Item(76, 13).Value = -2715.6924
$ws.Cells.Item(76, 14).Value = -3994.2856

$ws.Cells.Item(79, 8).Value = 3101.4546
$ws.Cells.Item(79, 9).Value = 3030.6924
$ws.Cells.Item(79, 10).Value = 3364.2856
$ws.Cells.Item(79, 11).Value = 3030.6924
$ws.Cells.Item(79, 12).Value = 3364.2856
$ws.Cells.Item(79, 13).Value = -1938.6924
$ws.Cells.Item(79, 14).Value = -5548.2856

$ws.Cells.Item(129, 8).Value = 725.8421
$ws.Cells.Item(129, 9).Value = 355.85715
$ws.Cells.Item(129, 10).Value = 941.6667
$ws.Cells.Item(129, 11).Value = 1067.57145
$ws.Cells.Item(129, 12).Value = 2825.0001
$ws.Cells.Item(129, 13).Value = 3932.42855
$ws.Cells.Item(129, 14).Value = -12825.0001

$ws.Cells.Item(137, 8).Value = 27028908
$ws.Cells.Item(137, 9).Value = 1146.5
$ws.Cells.Item(137, 11).Value = 3439.5
$ws.Cells.Item(137, 13).Value = -889.5

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(32, 8).Value = 1998.33
$ws.Cells.Item(32, 9).Value = 1962.3158
$ws.Cells.Item(32, 10).Value = 2682.6
$ws.Cells.Item(32, 11).Value = 1962.3158
$ws.Cells.Item(32, 12).Value = 2682.6
$ws.Cells.Item(32, 13).Value = -1675.3158
$ws.Cells.Item(32, 14).Value = -3256.6

$ws.Cells.Item(37, 8).Value = 9717.833000000001
$ws.Cells.Item(37, 9).Value = 0
$ws.Cells.Item(37, 10).Value = 9717.833000000001
$ws.Cells.Item(37, 11).Value = 0
$ws.Cells.Item(37, 12).Value = 9717.833000000001
$ws.Cells.Item(37, 13).ClearContents()
$ws.Cells.Item(37, 14).Value = -10263.833

$ws.Cells.Item(61, 8).Value = 1869.3043
$ws.Cells.Item(61, 9).Value = 1948.2142
$ws.Cells.Item(61, 10).Value = 1746.5555
$ws.Cells.Item(61, 11).Value = 1948.2142
$ws.Cells.Item(61, 12).Value = 1746.5555
$ws.Cells.Item(61, 13).Value = -1736.2142
$ws.Cells.Item(61, 14).Value = -2170.5555

$ws.Cells.Item(74, 8).Value = 4840.7026
$ws.Cells.Item(74, 9).Value = 893.96
$ws.Cells.Item(74, 11).Value = 893.96
$ws.Cells.Item(74, 13).Value = -19.96000000000004

$ws.Cells.Item(77, 8).Value = 4840.7026
$ws.Cells.Item(77, 9).Value = 893.96
$ws.Cells.Item(77, 11).Value = 4469.8
$ws.Cells.Item(77, 13).Value = -101.8000000000002

$ws.Cells.Item(122, 8).Value = 1840
$ws.Cells.Item(122, 9).Value = 1760
$ws.Cells.Item(122, 10).Value = 1920
$ws.Cells.Item(122, 11).Value = 5280
$ws.Cells.Item(122, 12).Value = 5760
$ws.Cells.Item(122, 13).Value = -2830
$ws.Cells.Item(122, 14).Value = -10660

$ws.Cells.Item(132, 8).Value = 20404.463
$ws.Cells.Item(132, 9).Value = 1827.0465
$ws.Cells.Item(132, 10).Value = 93025.27
$ws.Cells.Item(132, 11).Value = 5481.139499999999
$ws.Cells.Item(132, 12).Value = 279075.81
$ws.Cells.Item(132, 13).Value = -2951.139499999999
$ws.Cells.Item(132, 14).Value = -284135.81

$ws.Cells.Item(136, 8).Value = 1869.3043
$ws.Cells.Item(136, 9).Value = 1948.2142
$ws.Cells.Item(136, 10).Value = 1746.5555
$ws.Cells.Item(136, 11).Value = 5844.642599999999
$ws.Cells.Item(136, 12).Value = 5239.666499999999
$ws.Cells.Item(136, 13).Value = -3294.642599999999
$ws.Cells.Item(136, 14).Value = -10339.6665

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(31, 8).Value = 1767.5483
$ws.Cells.Item(31, 9).Value = 1246.5151
$ws.Cells.Item(31, 10).Value = 2360.4482
$ws.Cells.Item(31, 11).Value = 1246.5151
$ws.Cells.Item(31, 12).Value = 2360.4482
$ws.Cells.Item(31, 13).Value = -951.5151000000001
$ws.Cells.Item(31, 14).Value = -2950.4482

$ws.Cells.Item(34, 8).Value = 1767.5483
$ws.Cells.Item(34, 9).Value = 1246.5151
$ws.Cells.Item(34, 10).Value = 2360.4482
$ws.Cells.Item(34, 11).Value = 1246.5151
$ws.Cells.Item(34, 12).Value = 2360.4482
$ws.Cells.Item(34, 13).Value = -1044.5151
$ws.Cells.Item(34, 14).Value = -2764.4482

$ws.Cells.Item(51, 8).Value = 10722.357
$ws.Cells.Item(51, 10).Value = 11463.308
$ws.Cells.Item(51, 12).Value = 11463.308
$ws.Cells.Item(51, 14).Value = -12935.308

$ws.Cells.Item(59, 8).Value = 14461.909
$ws.Cells.Item(59, 10).Value = 15342.333
$ws.Cells.Item(59, 12).Value = 15342.333
$ws.Cells.Item(59, 14).Value = -17632.333

$ws.Cells.Item(61, 8).Value = 10722.357
$ws.Cells.Item(61, 10).Value = 11463.308
$ws.Cells.Item(61, 12).Value = 11463.308
$ws.Cells.Item(61, 14).Value = -12159.308

$ws.Cells.Item(74, 8).Value = 17932.7
$ws.Cells.Item(74, 9).Value = 4692.5
$ws.Cells.Item(74, 11).Value = 4692.5
$ws.Cells.Item(74, 13).Value = -3818.5

$ws.Cells.Item(77, 8).Value = 17932.7
$ws.Cells.Item(77, 9).Value = 4692.5
$ws.Cells.Item(77, 11).Value = 14077.5
$ws.Cells.Item(77, 13).Value = -9709.5

$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(113, 8).Value = 623.4783
$ws.Cells.Item(113, 10).Value = 593.125
$ws.Cells.Item(113, 12).Value = 1779.375
$ws.Cells.Item(113, 14).Value = -6119.375

$ws.Cells.Item(131, 8).Value = 2093.7446
$ws.Cells.Item(131, 9).Value = 12406
$ws.Cells.Item(131, 10).Value = 1514.4045
$ws.Cells.Item(131, 11).Value = 37218
$ws.Cells.Item(131, 12).Value = 4543.2135
$ws.Cells.Item(131, 13).Value = -32178
$ws.Cells.Item(131, 14).Value = -14623.2135

$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(40, 8).Value = 5000
$ws.Cells.Item(40, 10).Value = 5000
$ws.Cells.Item(40, 12).Value = 5000
$ws.Cells.Item(40, 14).Value = -5302

$ws.Cells.Item(70, 8).Value = 9485.944
$ws.Cells.Item(70, 9).Value = 10853.357
$ws.Cells.Item(70, 10).Value = 4700
$ws.Cells.Item(70, 11).Value = 10853.357
$ws.Cells.Item(70, 12).Value = 4700
$ws.Cells.Item(70, 13).Value = -10583.357
$ws.Cells.Item(70, 14).Value = -5240

$ws.Cells.Item(73, 8).Value = 9485.944
$ws.Cells.Item(73, 9).Value = 10853.357
$ws.Cells.Item(73, 10).Value = 4700
$ws.Cells.Item(73, 11).Value = 10853.357
$ws.Cells.Item(73, 12).Value = 4700
$ws.Cells.Item(73, 13).Value = -9917.357
$ws.Cells.Item(73, 14).Value = -6572

$ws.Cells.Item(113, 8).Value = 1146.5834
$ws.Cells.Item(113, 9).Value = 1149.7142
$ws.Cells.Item(113, 10).Value = 1142.2
$ws.Cells.Item(113, 11).Value = 1149.7142
$ws.Cells.Item(113, 12).Value = 1142.2
$ws.Cells.Item(113, 13).Value = 1020.2858
$ws.Cells.Item(113, 14).Value = -5482.2

$ws.Cells.Item(133, 8).Value = 52780
$ws.Cells.Item(133, 10).Value = 52780
$ws.Cells.Item(133, 12).Value = 52780
$ws.Cells.Item(133, 14).Value = -62900

$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(93, 8).Value = 1805.875
$ws.Cells.Item(93, 9).Value = 1322.2222
$ws.Cells.Item(93, 10).Value = 2427.7144
$ws.Cells.Item(93, 11).Value = 1322.2222
$ws.Cells.Item(93, 12).Value = 2427.7144
$ws.Cells.Item(93, 13).Value = -74.22219999999993
$ws.Cells.Item(93, 14).Value = -4923.7144

$ws.Cells.Item(122, 8).Value = 2840
$ws.Cells.Item(122, 9).Value = 2333.3333
$ws.Cells.Item(122, 10).Value = 3057.1428
$ws.Cells.Item(122, 11).Value = 6999.999899999999
$ws.Cells.Item(122, 12).Value = 9171.428400000001
$ws.Cells.Item(122, 13).Value = -4549.999899999999
$ws.Cells.Item(122, 14).Value = -14071.4284

$ws.Cells.Item(132, 8).Value = 5324.256
$ws.Cells.Item(132, 9).Value = 5515.1387
$ws.Cells.Item(132, 11).Value = 16545.4161
$ws.Cells.Item(132, 13).Value = -14015.4161

$ws.Cells.Item(136, 8).Value = 1788.1613
$ws.Cells.Item(136, 9).Value = 1082.5238
$ws.Cells.Item(136, 10).Value = 3270
$ws.Cells.Item(136, 11).Value = 3247.5714
$ws.Cells.Item(136, 12).Value = 9810
$ws.Cells.Item(136, 13).Value = -697.5713999999998
$ws.Cells.Item(136, 14).Value = -14910

$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(87, 8).Value = 43333.332
$ws.Cells.Item(87, 10).Value = 43333.332
$ws.Cells.Item(87, 12).Value = 43333.332
$ws.Cells.Item(87, 14).Value = -45829.332

$ws.Cells.Item(90, 8).Value = 43333.332
$ws.Cells.Item(90, 10).Value = 43333.332
$ws.Cells.Item(90, 12).Value = 129999.996
$ws.Cells.Item(90, 14).Value = -142479.996

$ws.Cells.Item(107, 8).Value = 617.4
$ws.Cells.Item(107, 9).Value = 435.41666
$ws.Cells.Item(107, 10).Value = 890.375
$ws.Cells.Item(107, 11).Value = 1306.24998
$ws.Cells.Item(107, 12).Value = 2671.125
$ws.Cells.Item(107, 13).Value = 613.7500199999999
$ws.Cells.Item(107, 14).Value = -6511.125

$ws.Cells.Item(136, 8).Value = 6645.864
$ws.Cells.Item(136, 9).Value = 7058.2896
$ws.Cells.Item(136, 10).Value = 4033.8333
$ws.Cells.Item(136, 11).Value = 21174.8688
$ws.Cells.Item(136, 12).Value = 12101.4999
$ws.Cells.Item(136, 13).Value = -18624.8688
$ws.Cells.Item(136, 14).Value = -17201.4999
